# Insert a new weekly price record as row 84, shifting existing rows
# (84-126) down to (85-127), per the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 84 (and everything below it) down by one row.
$ws.Rows(84).Insert()

# Populate the newly inserted row 84 with the new weekly record.
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C84").Value = "Los Lagos"
$ws.Range("D84").Value = 44466
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = 100112032
$ws.Range("G84").Value = "Zapallo italiano"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 17000
$ws.Range("L84").Value = 17000
$ws.Range("M84").Value = 17000
$ws.Range("N84").Value = "`$/caja 50 unidades"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 340
$ws.Range("Q84").Value = 50
$ws.Range("R84").Value = "Hortaliza"
